$d = $word.ActiveDocument

# Locate the exact run of text that needs to be split so Word's proofer
# markers (<w:proofErr>) can bracket the surname "Siemes", mirroring what
# Word itself does after a spell-check pass flags the name as unknown.
$old = "Weekopdracht Week 9 (Jort Siemes, s4028198)"

$target = $d.Content
$found = $target.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the title text to update"
}

# $target now spans exactly the original run's text (no paragraph mark),
# so replacing its contents via InsertXML only rewrites the runs inside
# the title paragraph - the paragraph's own pPr/style/ids are untouched.
$openXmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$wDocOpen = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$body = '<w:body><w:p>' `
    + '<w:r w:rsidRPr="00B360F0"><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Weekopdracht Week 9 (Jort </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Siemes</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>, s4028198)</w:t></w:r>' `
    + '</w:p></w:body>'
$wDocClose = '</w:document>'
$pkgClose = '</pkg:xmlData></pkg:part></pkg:package>'

$xml = $openXmlHeader + $pkgOpen + $wDocOpen + $body + $wDocClose + $pkgClose

$target.InsertXML($xml)
